$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Fix the "HMCTS SCSS" typo -> "HMCTS SSCS" (swap the transposed
#    "CS" -> "SC"), ending up split across three runs ("HMCTS SS" / "C"
#    / "S") exactly as in the target revision.
# ---------------------------------------------------------------------
$findRange = $d.Content
$found = $findRange.Find.Execute("HMCTS SCSS", $false, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'HMCTS SCSS' in the document"
}
$base = $findRange.Start

# Swap the transposed letters: "HMCTS S[C][S]S" -> "HMCTS S[S][C]S"
$rC = $d.Range($base + 7, $base + 8)
$rC.Text = "S"
$rS = $d.Range($base + 8, $base + 9)
$rS.Text = "C"

# Force the run to split at the two new boundaries (toggling a formatting
# property on and back off splits the run in two without altering its
# look, matching the three-run layout of the target revision).
$split1 = $d.Range($base + 8, $base + 9)
$split1.Font.Bold = 1
$split1b = $d.Range($base + 8, $base + 9)
$split1b.Font.Bold = 0

$split2 = $d.Range($base + 9, $base + 10)
$split2.Font.Bold = 1
$split2b = $d.Range($base + 9, $base + 10)
$split2b.Font.Bold = 0

# ---------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the "CM20 9RT" paragraph to the
#    end of the "PO BOX 12879" paragraph (right after "12879").
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $old = $d.Bookmarks.Item("_GoBack")
    $old.Delete()
}

$poBoxRange = $d.Content
$poBoxFound = $poBoxRange.Find.Execute("PO BOX 12879", $false, $false, $false, $false, `
                                        $false, $true, 1, $false, "", 0)
if (-not $poBoxFound) {
    throw "Could not find 'PO BOX 12879' in the document"
}
$endPos = $poBoxRange.End

# Insert a temporary marker character right after "12879", bookmark it,
# then delete the marker -- this leaves a zero-length "_GoBack" bookmark
# sitting exactly after "12879" and before the paragraph mark, without
# disturbing the existing "12879" run.
$ins = $d.Range($endPos, $endPos)
$ins.InsertAfter("X")
$markRange = $d.Range($endPos, $endPos + 1)
$d.Bookmarks.Add("_GoBack", $markRange)
$cleanup = $d.Range($endPos, $endPos + 1)
$cleanup.Text = ""
